$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Make room for the new "2022-Q1" sheet while keeping the sheetId /
#    rId bookkeeping identical to the target workbook: duplicate the
#    existing "总计" sheet (placed right after itself), rename the
#    original to "2022-Q1" and the duplicate back to "总计". This way
#    "2022-Q1" inherits the old "总计" sheetId and the new "总计" copy
#    gets the next one, exactly mirroring the authors' edit.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Copy($null, $total)
$totalCopy = $wb.Worksheets.Item("总计 (2)")

$total.Name = "2022-Q1"
$totalCopy.Name = "总计"

# ---------------------------------------------------------------------
# 2. Rebuild "2022-Q1" as a fund-holding sheet (same layout used by the
#    other quarterly sheets, e.g. "2021-Q4"): clear the old totals
#    content, copy the header/row formatting from "2021-Q4", then fill
#    in the new fund row.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Cells.Clear()

$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("A1:H2").Copy($q1.Range("A1"))
$q1.Range("A1").ClearContents()

$q1.Range("B2").Value = "'519097"
$q1.Range("C2").Value = "新华中小市值优选混合"
$q1.Range("D2").Value = "'0.75"
$q1.Range("E2").Value = "'62.70"
$q1.Range("F2").Value = "'4.27"
$q1.Range("G2").Value = "'0.0320"
$q1.Range("H2").Value = 5

# The leading apostrophes above force the numeric-looking strings to be
# stored as text (matching "519097", "0.75", "62.70", ... being plain
# text in the source data) without changing their value; drop the
# resulting "quote prefix" cell format so the cells fall back to the
# plain/default style, same as their counterparts on the other
# quarterly sheets.
$q1.Range("B2:G2").ClearFormats()

# ---------------------------------------------------------------------
# 3. Update "总计": insert a new row right after the header for the
#    2022-Q1 summary, pushing every other quarter down by one row, and
#    renumber the leading index column to match.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.03

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# Restore the originally active sheet/selection so we don't leave an
# unrelated "view state" diff behind.
[void]$wb.Worksheets.Item(1).Select()
[void]$wb.Worksheets.Item(1).Range("A1").Select()

